$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 39
$ws.Range("A39").Value2 = 111880591
$ws.Range("B39").Value2 = 90658
$ws.Range("E39").Value2 = 4361
$ws.Range("F39").Value2 = "Orange taggsvamp"
$ws.Range("G39").Value2 = "Hydnellum aurantiacum"
$ws.Range("H39").Value2 = "(Batsch:Fr.) P.Karst."
$ws.Range("I39").NumberFormat = "@"
$ws.Range("I39").Value2 = "8"
$ws.Range("Q39").Value2 = 509822.1902239832
$ws.Range("R39").Value2 = 6753234.069152902
$ws.Range("AJ39").Value2 = "tall"
$ws.Range("AK39").Value2 = "Pinus sylvestris"
$ws.Range("AO39").Value2 = "Pinus sylvestris"

# Row 40
$ws.Range("A40").Value2 = 111880509
$ws.Range("B40").Value2 = 90652
$ws.Range("E40").Value2 = 3100
$ws.Range("F40").Value2 = "Talltaggsvamp"
$ws.Range("G40").Value2 = "Bankera fuligineoalba"
$ws.Range("H40").Value2 = "(Schmidt : Fr.) Pouzar"
$ws.Range("I40").NumberFormat = "@"
$ws.Range("I40").Value2 = "6"
$ws.Range("Q40").Value2 = 509834.2096935506
$ws.Range("R40").Value2 = 6753644.114383955
$ws.Range("AJ40").Value2 = "tall"
$ws.Range("AK40").Value2 = "Pinus sylvestris"
$ws.Range("AO40").Value2 = "Pinus sylvestris"

# Row 41
$ws.Range("A41").Value2 = 111880532
$ws.Range("B41").Value2 = 90652
$ws.Range("E41").Value2 = 3100
$ws.Range("F41").Value2 = "Talltaggsvamp"
$ws.Range("G41").Value2 = "Bankera fuligineoalba"
$ws.Range("H41").Value2 = "(Schmidt : Fr.) Pouzar"
$ws.Range("I41").NumberFormat = "@"
$ws.Range("I41").Value2 = "2"
$ws.Range("Q41").Value2 = 509682.5105515064
$ws.Range("R41").Value2 = 6753540.591470475
$ws.Range("AJ41").Value2 = "tall"
$ws.Range("AK41").Value2 = "Pinus sylvestris"
$ws.Range("AO41").Value2 = "Pinus sylvestris"

# Row 42
$ws.Range("A42").Value2 = 111880580
$ws.Range("B42").Value2 = 90658
$ws.Range("E42").Value2 = 4361
$ws.Range("F42").Value2 = "Orange taggsvamp"
$ws.Range("G42").Value2 = "Hydnellum aurantiacum"
$ws.Range("H42").Value2 = "(Batsch:Fr.) P.Karst."
$ws.Range("Q42").Value2 = 509755.441071702
$ws.Range("R42").Value2 = 6753236.317390828
$ws.Range("AJ42").Value2 = "tall"
$ws.Range("AK42").Value2 = "Pinus sylvestris"
$ws.Range("AO42").Value2 = "Pinus sylvestris"

# Row 43
$ws.Range("A43").Value2 = 111880484
$ws.Range("B43").Value2 = 90658
$ws.Range("E43").Value2 = 4361
$ws.Range("F43").Value2 = "Orange taggsvamp"
$ws.Range("G43").Value2 = "Hydnellum aurantiacum"
$ws.Range("H43").Value2 = "(Batsch:Fr.) P.Karst."
$ws.Range("I43").NumberFormat = "@"
$ws.Range("I43").Value2 = "11"
$ws.Range("Q43").Value2 = 509900.7891887496
$ws.Range("R43").Value2 = 6753525.142772059
$ws.Range("AJ43").Value2 = "tall"
$ws.Range("AK43").Value2 = "Pinus sylvestris"
$ws.Range("AO43").Value2 = "Pinus sylvestris"

# Row 44
$ws.Range("A44").Value2 = 111880574
$ws.Range("B44").Value2 = 90658
$ws.Range("E44").Value2 = 4361
$ws.Range("F44").Value2 = "Orange taggsvamp"
$ws.Range("G44").Value2 = "Hydnellum aurantiacum"
$ws.Range("H44").Value2 = "(Batsch:Fr.) P.Karst."
$ws.Range("I44").NumberFormat = "@"
$ws.Range("I44").Value2 = "2"
$ws.Range("Q44").Value2 = 509595.7160662179
$ws.Range("R44").Value2 = 6753391.52735021
$ws.Range("AJ44").Value2 = "tall"
$ws.Range("AK44").Value2 = "Pinus sylvestris"
$ws.Range("AO44").Value2 = "Pinus sylvestris"

# Row 45
$ws.Range("A45").Value2 = 111880475
$ws.Range("B45").Value2 = 88966
$ws.Range("E45").Value2 = 5754
$ws.Range("F45").Value2 = "Gultoppig fingersvamp"
$ws.Range("G45").Value2 = "Ramaria testaceoflava"
$ws.Range("H45").Value2 = "(Bres.) Corner"
$ws.Range("Q45").Value2 = 509957.7514087428
$ws.Range("R45").Value2 = 6753362.853637428
$ws.Range("AJ45").Value2 = "gran"
$ws.Range("AK45").Value2 = "Picea abies"
$ws.Range("AO45").Value2 = "Picea abies"

# Row 46
$ws.Range("A46").Value2 = 111880562
$ws.Range("B46").Value2 = 90658
$ws.Range("E46").Value2 = 4361
$ws.Range("F46").Value2 = "Orange taggsvamp"
$ws.Range("G46").Value2 = "Hydnellum aurantiacum"
$ws.Range("H46").Value2 = "(Batsch:Fr.) P.Karst."
$ws.Range("I46").NumberFormat = "@"
$ws.Range("I46").Value2 = "3"
$ws.Range("Q46").Value2 = 509657.7198006394
$ws.Range("R46").Value2 = 6753521.069647122
$ws.Range("AJ46").Value2 = "tall"
$ws.Range("AK46").Value2 = "Pinus sylvestris"
$ws.Range("AO46").Value2 = "Pinus sylvestris"

# Row 47
$ws.Range("A47").Value2 = 111880601
$ws.Range("B47").Value2 = 88966
$ws.Range("E47").Value2 = 5754
$ws.Range("F47").Value2 = "Gultoppig fingersvamp"
$ws.Range("G47").Value2 = "Ramaria testaceoflava"
$ws.Range("H47").Value2 = "(Bres.) Corner"
$ws.Range("I47").NumberFormat = "@"
$ws.Range("I47").Value2 = "4"
$ws.Range("Q47").Value2 = 509941.5744066621
$ws.Range("R47").Value2 = 6753224.672924293
$ws.Range("AJ47").Value2 = "tall"
$ws.Range("AK47").Value2 = "Pinus sylvestris"
$ws.Range("AO47").Value2 = "Pinus sylvestris"

# Row 48
$ws.Range("A48").Value2 = 111880462
$ws.Range("B48").Value2 = 88966
$ws.Range("E48").Value2 = 5754
$ws.Range("F48").Value2 = "Gultoppig fingersvamp"
$ws.Range("G48").Value2 = "Ramaria testaceoflava"
$ws.Range("H48").Value2 = "(Bres.) Corner"
$ws.Range("I48").NumberFormat = "@"
$ws.Range("I48").Value2 = "1"
$ws.Range("Q48").Value2 = 509970.2466718731
$ws.Range("R48").Value2 = 6753250.046013334
$ws.Range("AJ48").Value2 = "tall"
$ws.Range("AK48").Value2 = "Pinus sylvestris"
$ws.Range("AO48").Value2 = "Pinus sylvestris # vid tallar"

# Row 49
$ws.Range("A49").Value2 = 111880500
$ws.Range("B49").Value2 = 88966
$ws.Range("E49").Value2 = 5754
$ws.Range("F49").Value2 = "Gultoppig fingersvamp"
$ws.Range("G49").Value2 = "Ramaria testaceoflava"
$ws.Range("H49").Value2 = "(Bres.) Corner"
$ws.Range("I49").NumberFormat = "@"
$ws.Range("I49").Value2 = "4"
$ws.Range("Q49").Value2 = 509899.1991435916
$ws.Range("R49").Value2 = 6753571.34232254
$ws.Range("AJ49").Value2 = "gran"
$ws.Range("AK49").Value2 = "Picea abies"
$ws.Range("AO49").Value2 = "Picea abies"

# AL44 -> AL48 move ("vid tallar" note)
$ws.Range("AL44").ClearContents()
$ws.Range("AL48").Value2 = "vid tallar"
